# devedores.xlsx — "Inicializando parte 2 do projeto" commit
#
# Row 4 (Kiko) gains a long run of (date, amount) payment-history pairs in
# columns E..Z, and the running-total cell A4 is bumped from 250 to 550.0.
# A brand-new debtor row (row 6, Hugo) is appended below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: update running total, then append the E4:Z4 history pairs ----
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "550.0"
$ws.Range("A4").ClearFormats()

$dateCols = @("E","G","I","K","M","O","Q","S","U","W","Y")
$amtCols  = @("F","H","J","L","N","P","R","T","V","X","Z")

for ($i = 0; $i -lt $dateCols.Length; $i++) {
    $dCell = $ws.Range($dateCols[$i] + "4")
    $dCell.Value = "24/12/2020"

    $aCell = $ws.Range($amtCols[$i] + "4")
    $aCell.NumberFormat = "@"
    $aCell.Value = "300.0"
    $aCell.ClearFormats()
}

# --- Row 6: brand-new debtor entry (Hugo) ---------------------------------
$ws.Range("A6").Value = 500
$ws.Range("B6").Value = "Hugo"
$ws.Range("C6").Value = "24/12/2020"
$ws.Range("D6").Value = 500
